# Apply cryptos list update (prices/volumes refreshed; Hedera/TrustWalletToken rows swapped)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.491.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.574.46"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "292.28"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3723"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.91"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3408"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.151"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07562"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.30"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.979"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.571.88"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06753"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.307"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.39"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.474.19"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.630"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.04"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.20"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.067"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.60"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.745.29"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.082"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +9.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.240"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.014"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.863"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08382"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02488"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2304"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.339"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.66%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06563"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.466"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.16%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6260"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.99"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.815"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5850"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.35"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.081"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.218"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07342"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.18%  "
